# Log a new day entry (row 10): "Working on camera colision and movement"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

$row = 10

# Copy the existing formatting from the row above so the new row picks up
# the same number formats / styles (date, time) instead of minting new ones.
$ws.Cells.Item($row - 1, 2).Copy() | Out-Null
$ws.Cells.Item($row, 2).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item($row - 1, 3).Copy() | Out-Null
$ws.Cells.Item($row, 3).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item($row - 1, 4).Copy() | Out-Null
$ws.Cells.Item($row, 4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item($row - 1, 5).Copy() | Out-Null
$ws.Cells.Item($row, 5).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item($row, 1).Value = "Thursday"
$ws.Cells.Item($row, 2).Value = 45778
$ws.Cells.Item($row, 3).Value = 0.79166666666666663
$ws.Cells.Item($row, 4).Value = 0.84027777777777779
$ws.Cells.Item($row, 5).Formula = "=D10-C10"
$ws.Cells.Item($row, 6).Value = "Working on camera colision and movement"

$ws.Range("E10").Select()
